$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values from F2:M3 while keeping cell formatting/styles intact
$ws.Range("F2:M3").ClearContents()

# Update the selection to F2:M3 with active cell F2
$ws.Range("F2:M3").Select()
